$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.505.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.101.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.87"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5228"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4490"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.73"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +16.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08926"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.156"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.49"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.093.18"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.723"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.714"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.41"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06626"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.65%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.292"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.550.54"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.35"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.334"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.340.08"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.94%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.578"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.70%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.199"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1074"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.670"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.164"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.903"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.53"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02572"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06783"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.484"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.72"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2266"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6927"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.254"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.39%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.279"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.637"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.239"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.245"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.03"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.50%  "
